$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '36.447.93'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '1.944.52'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '243.14'
$ws.Range('E5').Value = '  +0.69%  '
Set-TextValue 'D6' '0.619'
$ws.Range('E6').Value = '  -0.76%  '
Set-TextValue 'D7' '58.50'
$ws.Range('E7').Value = '  -3.00%  '
$ws.Range('E8').Value = '  -0.03%  '
Set-TextValue 'D9' '0.365'
$ws.Range('E9').Value = '  -2.34%  '
Set-TextValue 'D10' '55.82'
$ws.Range('E10').Value = '  -1.15%  '
Set-TextValue 'D11' '0.0831'
$ws.Range('E11').Value = '  +3.55%  '
Set-TextValue 'D12' '0.103'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D13' '0.822'
$ws.Range('E13').Value = '  -3.86%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D14' '21.54'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('D15').Value = '2.227.52'
$ws.Range('E15').Value = '  -0.79%  '
Set-TextValue 'D16' '13.61'
$ws.Range('E16').Value = '  -2.78%  '
Set-TextValue 'D17' '5.24'
$ws.Range('E17').Value = '  -3.04%  '
$ws.Range('D18').Value = '1.946.72'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = '36.324.57'
$ws.Range('E19').Value = '  +1.04%  '
Set-TextValue 'D20' '69.52'
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').Value = '0.0₃0861'
$ws.Range('E21').Value = '  +0.48%  '
Set-TextValue 'D22' '228.52'
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -3.37%  '
Set-TextValue 'D26' '2.29'
$ws.Range('E26').Value = '  +0.22%  '
Set-TextValue 'D27' '9.20'
$ws.Range('E27').Value = '  -5.60%  '
Set-TextValue 'D28' '161.66'
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D29' '0.130'
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D30' '19.47'
$ws.Range('E30').Value = '  -1.32%  '
Set-TextValue 'D31' '0.118'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('E32').Value = '  +1.40%  '
Set-TextValue 'D33' '4.67'
$ws.Range('E33').Value = '  -3.97%  '
Set-TextValue 'D34' '0.0627'
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  -2.25%  '
$ws.Range('E39').Value = '  -6.05%  '
Set-TextValue 'D40' '3.03'
$ws.Range('E40').Value = '  -1.33%  '
Set-TextValue 'D41' '0.0982'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  +0.82%  '
$ws.Range('E43').Value = '  -4.25%  '
$ws.Range('E44').Value = '  -1.43%  '
Set-TextValue 'D45' '16.01'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').Value = '1.357.56'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('E47').Value = '  -4.68%  '
Set-TextValue 'D48' '87.84'
$ws.Range('E48').Value = '  -4.44%  '
Set-TextValue 'D49' '7.12'
$ws.Range('E49').Value = '  -4.50%  '
Set-TextValue 'D50' '2.82'
$ws.Range('E50').Value = '  -0.52%  '
Set-TextValue 'D51' '45.41'
$ws.Range('E51').Value = '  +3.61%  '
